# COREESG_holdings.xlsx refresh
#  - Bump the "Model holdings provided as of" date in the confidential
#    disclosure footer from 2021-04-05 to 2021-04-06.
#  - Refresh the Weight / Percent Change figures for each ETF sleeve
#    (rows 2-7 of Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet carries protection, so lift it for the duration of the edits
# and restore it once we're done.
$ws.Unprotect()

# --- Update the confidential disclosure footer date -----------------------
$oldText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-05 for illustrative purposes only and are subject to change."
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-06 for illustrative purposes only and are subject to change."

foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value2
        if ($val -ne $null -and $val.ToString() -eq $oldText) {
            $cell.Value = $newText
        }
    }
}

# --- Refresh Weight / Percent Change values --------------------------------
$ws.Range("D2").Value = 0.2509894946380719
$ws.Range("E2").Value = 0.0006676681689201658

$ws.Range("D3").Value = 0.4930196760685087
$ws.Range("E3").Value = -0.001902173913043592

$ws.Range("D4").Value = 0.09990584440589723
$ws.Range("E4").Value = 0.00640279394644927

$ws.Range("D5").Value = 0.098763073209841
$ws.Range("E5").Value = 0.002858776443682043

$ws.Range("D6").Value = 0.05732191167768132
$ws.Range("E6").Value = 0.0002279462046959058

$ws.Range("E7").Value = 0.0001648529250564135

# Restore the original sheet protection.
$ws.Protect()
